# Auto-generated edit script applying scheduled-runner price/profit updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 25311.54
$ws.Range("J93").Value = 25311.54
$ws.Range("L93").Value = 25311.54
$ws.Range("N93").Value = -30303.54
$ws.Range("H98").Value = 4175
$ws.Range("I98").Value = 1719.4445
$ws.Range("J98").Value = 6016.6665
$ws.Range("K98").Value = 1719.4445
$ws.Range("L98").Value = 6016.6665
$ws.Range("M98").Value = -221.4445000000001
$ws.Range("N98").Value = -9012.666499999999
$ws.Range("H112").Value = 1283.9016
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1291.9667
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 3875.9001
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6091.9001
$ws.Range("H113").Value = 5276.4375
$ws.Range("I113").Value = 2492.4546
$ws.Range("J113").Value = 11401.2
$ws.Range("K113").Value = 2492.4546
$ws.Range("L113").Value = 11401.2
$ws.Range("M113").Value = 761.5454
$ws.Range("N113").Value = -17909.2
$ws.Range("H122").Value = 4175
$ws.Range("I122").Value = 1719.4445
$ws.Range("J122").Value = 6016.6665
$ws.Range("K122").Value = 5158.333500000001
$ws.Range("L122").Value = 18049.9995
$ws.Range("M122").Value = -2708.333500000001
$ws.Range("N122").Value = -22949.9995
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H129").Value = 822.3838500000001
$ws.Range("I129").Value = 311
$ws.Range("J129").Value = 867.34064
$ws.Range("K129").Value = 933
$ws.Range("L129").Value = 2602.02192
$ws.Range("M129").Value = 4067
$ws.Range("N129").Value = -12602.02192
$ws.Range("H131").Value = 4580
$ws.Range("I131").Value = 4250
$ws.Range("J131").Value = 4800
$ws.Range("K131").Value = 12750
$ws.Range("L131").Value = 14400
$ws.Range("M131").Value = -7710
$ws.Range("N131").Value = -24480
$ws.Range("H133").Value = 41338.75
$ws.Range("J133").Value = 41338.75
$ws.Range("L133").Value = 41338.75
$ws.Range("N133").Value = -51458.75
$ws.Range("H135").Value = 1056.4762
$ws.Range("I135").Value = 577.5714
$ws.Range("K135").Value = 5198.1426
$ws.Range("M135").Value = -2663.1426
$ws.Range("H136").Value = 51120
$ws.Range("J136").Value = 51120
$ws.Range("L136").Value = 51120
$ws.Range("N136").Value = -61320
$ws.Range("H140").Value = 47648.57
$ws.Range("J140").Value = 47648.57
$ws.Range("L140").Value = 47648.57
$ws.Range("N140").Value = -58008.57
$ws.Range("H141").Value = 169090.58
$ws.Range("I141").Value = 223999.11
$ws.Range("J141").Value = 4365
$ws.Range("K141").Value = 671997.33
$ws.Range("L141").Value = 13095
$ws.Range("M141").Value = -666817.33
$ws.Range("N141").Value = -23455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2024.3334
$ws.Range("I74").Value = 774.1111
$ws.Range("J74").Value = 3899.6667
$ws.Range("K74").Value = 774.1111
$ws.Range("L74").Value = 3899.6667
$ws.Range("M74").Value = 99.88890000000004
$ws.Range("N74").Value = -5647.6667
$ws.Range("H77").Value = 2024.3334
$ws.Range("I77").Value = 774.1111
$ws.Range("J77").Value = 3899.6667
$ws.Range("K77").Value = 3870.5555
$ws.Range("L77").Value = 19498.3335
$ws.Range("M77").Value = 497.4445000000001
$ws.Range("N77").Value = -28234.3335
$ws.Range("H92").Value = 28937.5
$ws.Range("J92").Value = 28937.5
$ws.Range("L92").Value = 28937.5
$ws.Range("N92").Value = -33929.5
$ws.Range("H132").Value = 4445.857
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4445.857
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 13337.571
$ws.Range("N132").Value = -18397.571
$ws.Range("H137").Value = 38570
$ws.Range("J137").Value = 40762.5
$ws.Range("L137").Value = 40762.5
$ws.Range("N137").Value = -50962.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 543.46155
$ws.Range("J64").Value = 505.85715
$ws.Range("L64").Value = 505.85715
$ws.Range("N64").Value = -955.85715
$ws.Range("H67").Value = 543.46155
$ws.Range("J67").Value = 505.85715
$ws.Range("L67").Value = 505.85715
$ws.Range("N67").Value = -2065.85715
$ws.Range("H137").Value = 45500
$ws.Range("J137").Value = 45500
$ws.Range("L137").Value = 45500
$ws.Range("N137").Value = -55700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3348.818
$ws.Range("I31").Value = 1540.2858
$ws.Range("J31").Value = 6513.75
$ws.Range("K31").Value = 1540.2858
$ws.Range("L31").Value = 6513.75
$ws.Range("M31").Value = -1245.2858
$ws.Range("N31").Value = -7103.75
$ws.Range("H34").Value = 3348.818
$ws.Range("I34").Value = 1540.2858
$ws.Range("J34").Value = 6513.75
$ws.Range("K34").Value = 1540.2858
$ws.Range("L34").Value = 6513.75
$ws.Range("M34").Value = -1338.2858
$ws.Range("N34").Value = -6917.75
$ws.Range("H41").Value = 31421.143
$ws.Range("I41").Value = 9722
$ws.Range("J41").Value = 40100.8
$ws.Range("K41").Value = 9722
$ws.Range("L41").Value = 40100.8
$ws.Range("M41").Value = -9294
$ws.Range("N41").Value = -40956.8
$ws.Range("H99").Value = 16672292
$ws.Range("I99").Value = 50002000
$ws.Range("J99").Value = 7437.375
$ws.Range("K99").Value = 50002000
$ws.Range("L99").Value = 7437.375
$ws.Range("M99").Value = -50000502
$ws.Range("N99").Value = -10433.375
$ws.Range("H122").Value = 4280
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -49900
$ws.Range("H126").Value = 16672292
$ws.Range("I126").Value = 50002000
$ws.Range("J126").Value = 7437.375
$ws.Range("K126").Value = 150006000
$ws.Range("L126").Value = 22312.125
$ws.Range("M126").Value = -150003530
$ws.Range("N126").Value = -27252.125
$ws.Range("H132").Value = 4031.28
$ws.Range("I132").Value = 3353.818
$ws.Range("K132").Value = 10061.454
$ws.Range("M132").Value = -7531.454000000002
$ws.Range("H134").Value = 9945.357
$ws.Range("I134").Value = 11573.5
$ws.Range("J134").Value = 5875
$ws.Range("K134").Value = 34720.5
$ws.Range("L134").Value = 17625
$ws.Range("M134").Value = -32185.5
$ws.Range("N134").Value = -22695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 514510.94
$ws.Range("I5").Value = 532.4
$ws.Range("K5").Value = 1597.2
$ws.Range("M5").Value = -1485.2
$ws.Range("H63").Value = 4975
$ws.Range("I63").Value = 5033.3335
$ws.Range("J63").Value = 4940
$ws.Range("K63").Value = 15100.0005
$ws.Range("L63").Value = 14820
$ws.Range("M63").Value = -14351.0005
$ws.Range("N63").Value = -16318
$ws.Range("H66").Value = 4975
$ws.Range("I66").Value = 5033.3335
$ws.Range("J66").Value = 4940
$ws.Range("K66").Value = 45300.0015
$ws.Range("L66").Value = 44460
$ws.Range("M66").Value = -41556.0015
$ws.Range("N66").Value = -51948
$ws.Range("H122").Value = 3007.9583
$ws.Range("I122").Value = 982.2222
$ws.Range("J122").Value = 3475.4358
$ws.Range("K122").Value = 8839.9998
$ws.Range("L122").Value = 31278.9222
$ws.Range("M122").Value = -6389.9998
$ws.Range("N122").Value = -36178.9222
$ws.Range("H135").Value = 514510.94
$ws.Range("I135").Value = 532.4
$ws.Range("K135").Value = 4791.599999999999
$ws.Range("M135").Value = -2256.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 5240
$ws.Range("J17").Value = 5240
$ws.Range("L17").Value = 5240
$ws.Range("N17").Value = -5576
$ws.Range("H102").Value = 2874.5
$ws.Range("I102").Value = 1899.875
$ws.Range("J102").Value = 6773
$ws.Range("K102").Value = 1899.875
$ws.Range("L102").Value = 6773
$ws.Range("M102").Value = -277.875
$ws.Range("N102").Value = -10017
$ws.Range("H107").Value = 9259961
$ws.Range("I107").Value = 330
$ws.Range("J107").Value = 12346505
$ws.Range("K107").Value = 330
$ws.Range("L107").Value = 12346505
$ws.Range("M107").Value = 1590
$ws.Range("N107").Value = -12350345
$ws.Range("H122").Value = 4441.6665
$ws.Range("I122").Value = 1610.6364
$ws.Range("J122").Value = 12227
$ws.Range("K122").Value = 4831.9092
$ws.Range("L122").Value = 36681
$ws.Range("M122").Value = -2381.9092
$ws.Range("N122").Value = -41581
$ws.Range("H137").Value = 40456
$ws.Range("J137").Value = 40456
$ws.Range("L137").Value = 40456
$ws.Range("N137").Value = -50656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 12004.333
$ws.Range("I32").Value = 12004.333
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 12004.333
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -11687.333
$ws.Range("H40").Value = 5280.864
$ws.Range("I40").Value = 4429.231
$ws.Range("K40").Value = 4429.231
$ws.Range("M40").Value = -4293.231
$ws.Range("H46").Value = 1637
$ws.Range("I46").Value = 1115.8334
$ws.Range("J46").Value = 2262.4
$ws.Range("K46").Value = 1115.8334
$ws.Range("L46").Value = 2262.4
$ws.Range("M46").Value = -927.8334
$ws.Range("N46").Value = -2638.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 20000
$ws.Range("N74").Value = -21872
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 60000
$ws.Range("N77").Value = -69360
$ws.Range("H113").Value = 8116.4614
$ws.Range("J113").Value = 1161
$ws.Range("L113").Value = 3483
$ws.Range("N113").Value = -7823
$ws.Range("H122").Value = 6870.3
$ws.Range("I122").Value = 5100.1665
$ws.Range("K122").Value = 15300.4995
$ws.Range("M122").Value = -12850.4995
